$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("E24").Value = 15.3
$ws.Range("E26").Value = 22.42
